$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "function_parameters": update parameter types (ENG/String -> int /
# float / DataFrame / string as appropriate) and append a new row (id 8)
# recording a float parameter.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("function_parameters")

# Row 2: Octopus_Params (function 1) -> type becomes int
$ws3.Range("E2").Value = "int"

# Row 3: Sys_Params (function 1) -> type becomes int
$ws3.Range("E3").Value = "int"

# Row 4: text (function 2) -> value becomes numeric 4, type becomes int
$ws3.Range("D4").Value = 4
$ws3.Range("E4").Value = "int"

# Row 5 & 6 keep type DataFrame (shared string index shifts only, no value
# change needed) and row 7 keeps type string - nothing to do for those.

# New row 8: function 2, Octopus_Params, value 5.5, type float.
# Copy formatting from row 7 first so the new row matches the table style.
$ws3.Range("B7:E7").Copy()
$ws3.Range("B8:E8").PasteSpecial(-4122)  # xlPasteFormats
$ws3.Range("B8").Value = 2
$ws3.Range("C8").Value = "Octopus_Params"
$ws3.Range("D8").Value = 5.5
$ws3.Range("E8").Value = "float"

# ---------------------------------------------------------------------------
# View / selection state: Functions sheet loses the active-tab flag and its
# scrolled-right view, selection moves to C3; function_parameters becomes the
# active sheet with selection on D10 (one row below the newly-added data).
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Functions")
[void]$ws1.Activate()
[void]$ws1.Range("C3").Select()

[void]$ws3.Activate()
[void]$ws3.Range("D10").Select()
